# Update countries & provincias Spain
# Applies the data refresh captured in the commit:
#  - timestamp of the data refresh (A1) bumped from 05:45 to 07:02
#  - updated case counters for several countries (India, Pakistan,
#    Uzbekistan, Tailandia, Mongolia, Butan)
#  - countries "Belice" and "Islas Malvinas" were re-sorted earlier in
#    the list (their shared-string position moved up, shifting the
#    rows that follow them down by one row), so the rows 167-170 and
#    213-214 are rewritten in full (name + stats) to reflect the new
#    row order and Belice's refreshed stats.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 19 de Agosto de 2020 a las 07:02"

# --- India (row 6) ------------------------------------------------------
$ws.Range("B6").Value = 2767273
$ws.Range("C6").Value = 647
$ws.Range("E6").Value = 677556

# --- Pakistan (row 18) ---------------------------------------------------
$ws.Range("B18").Value = 290445
$ws.Range("C18").Value = 613
$ws.Range("D18").Value = 272128
$ws.Range("E18").Value = 12116
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = 6201

# --- Uzbekistan (row 61) --------------------------------------------------
$ws.Range("B61").Value = 36518
$ws.Range("C61").Value = 166
$ws.Range("D61").Value = 32223
$ws.Range("E61").Value = 4052
$ws.Range("G61").Value = 1
$ws.Range("H61").Value = 243

# --- Tailandia (row 118) --------------------------------------------------
$ws.Range("B118").Value = 3382
$ws.Range("C118").Value = 1
$ws.Range("D118").Value = 3199

# --- Belice / Guadalupe / Tanzania / Taiwan block (rows 167-170) ---------
# Belice's shared string moved ahead of Guadalupe, so it now occupies row
# 167 with refreshed stats, and Guadalupe/Tanzania/Taiwan each shift down
# one row keeping their (unchanged) stats.
$ws.Range("A167").Value = "Belice"
$ws.Range("B167").Value = 553
$ws.Range("C167").Value = 78
$ws.Range("D167").Value = 38
$ws.Range("E167").Value = 511
$ws.Range("H167").Value = 4

$ws.Range("A168").Value = "Guadalupe"
$ws.Range("B168").Value = 510
$ws.Range("D168").Value = 289
$ws.Range("E168").Value = 206
$ws.Range("H168").Value = 15

$ws.Range("A169").Value = "Tanzania"
$ws.Range("B169").Value = 509
$ws.Range("D169").Value = 183
$ws.Range("E169").Value = 305
$ws.Range("H169").Value = 21

$ws.Range("A170").Value = "Taiwan"
$ws.Range("B170").Value = 486
$ws.Range("D170").Value = 450
$ws.Range("E170").Value = 29
$ws.Range("H170").Value = 7

# --- Mongolia (row 182) ---------------------------------------------------
$ws.Range("D182").Value = 281
$ws.Range("E182").Value = 17

# --- Butan (row 190) -------------------------------------------------------
$ws.Range("B190").Value = 147
$ws.Range("C190").Value = 1
$ws.Range("E190").Value = 44

# --- Islas Malvinas / Montserrat swap (rows 213-214) ----------------------
# Islas Malvinas' shared string moved ahead of Montserrat, so the two rows
# swap their country name while keeping their own (unchanged) stats.
$ws.Range("A213").Value = "Islas Malvinas"
$ws.Range("D213").Value = 13
$ws.Range("H213").Value = 0

$ws.Range("A214").Value = "Montserrat"
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1
